# Generate Report for Handoff
# Adds two new "Ready for handoff" rows (38efd1a5-... and 5e154c0c-...)
# ahead of the existing e011a692-... row on all three sheets
# (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Push the existing "e011a692" row (currently row 3) down to row 5,
# opening up rows 3 and 4 for the two new entries.
$ws.Range("A3:A4").EntireRow.Insert()

# Rebuild the hyperlinks collection (row-insert does not re-anchor the
# existing hyperlink ranges), so drop everything on this sheet first.
$ws.Hyperlinks.Delete()

$ws.Range("A3").Value = "38efd1a5-cded-43d5-9e40-a56e8c5d8d85.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-21 10:34:43"

$ws.Range("A4").Value = "5e154c0c-7e94-41fe-895d-735eefa2dd79.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "2016-03-21 10:34:43"

# Row 5 already holds the shifted-down e011a692 values/format, untouched.

$ws.Hyperlinks.Add($ws.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/007b4ece8461d1d54c433d65ab9d1e56c718bffb/e2e/38cdfdee-e750-404e-bd36-1a8a0262a62d.md", "", "", "38cdfdee-e750-404e-bd36-1a8a0262a62d.md")
$ws.Hyperlinks.Add($ws.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/38efd1a5-cded-43d5-9e40-a56e8c5d8d85.md", "", "", "38efd1a5-cded-43d5-9e40-a56e8c5d8d85.md")
$ws.Hyperlinks.Add($ws.Cells.Item(4,1), "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/5e154c0c-7e94-41fe-895d-735eefa2dd79.md", "", "", "5e154c0c-7e94-41fe-895d-735eefa2dd79.md")
$ws.Hyperlinks.Add($ws.Cells.Item(5,1), "https://github.com/OpenLocalizationTest/oltest/blob/7e65c197f1b142ae59144fc72bb407b8868b3994/e2e/e011a692-b796-4695-a3af-0d2a96468c8d.md", "", "", "e011a692-b796-4695-a3af-0d2a96468c8d.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A3:A4").EntireRow.Insert()
$ws.Hyperlinks.Delete()

$ws.Range("A3").Value = "38efd1a5-cded-43d5-9e40-a56e8c5d8d85.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "38efd1a5-cded-43d5-9e40-a56e8c5d8d85.328696bf5fa0985729bbcefda211a578ce768408.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-21 10:34:40"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("J3").Value = "Include"

$ws.Range("A4").Value = "5e154c0c-7e94-41fe-895d-735eefa2dd79.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "5e154c0c-7e94-41fe-895d-735eefa2dd79.9c7f514114b497b10602888709ccdf35bc06a240.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-21 10:34:40"
$ws.Range("H4").Value = "0001-01-01 00:00:00"
$ws.Range("J4").Value = "Include"

# Row 5 already holds the shifted-down e011a692 values/format, untouched.

$ws.Hyperlinks.Add($ws.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/007b4ece8461d1d54c433d65ab9d1e56c718bffb/e2e/38cdfdee-e750-404e-bd36-1a8a0262a62d.md", "", "", "38cdfdee-e750-404e-bd36-1a8a0262a62d.md")
$ws.Hyperlinks.Add($ws.Cells.Item(2,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/231d1f53aa916a3eccdce0a9244e97996e7fa4ff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/38cdfdee-e750-404e-bd36-1a8a0262a62d.e55059b9176280270322352189b79dcaa829d7ae.zh-cn.xlf", "", "", "38cdfdee-e750-404e-bd36-1a8a0262a62d.e55059b9176280270322352189b79dcaa829d7ae.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7716478d7d6f68e06f9142cdc7eab72f6aac80dd/e2e/38cdfdee-e750-404e-bd36-1a8a0262a62d.md", "", "", "38cdfdee-e750-404e-bd36-1a8a0262a62d.md")
$ws.Hyperlinks.Add($ws.Cells.Item(2,7), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6f8670408e94cb08e7bf013eb9a9b5cc6392a965/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/38cdfdee-e750-404e-bd36-1a8a0262a62d.e55059b9176280270322352189b79dcaa829d7ae.zh-cn.xlf", "", "", "38cdfdee-e750-404e-bd36-1a8a0262a62d.e55059b9176280270322352189b79dcaa829d7ae.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/38efd1a5-cded-43d5-9e40-a56e8c5d8d85.md", "", "", "38efd1a5-cded-43d5-9e40-a56e8c5d8d85.md")
$ws.Hyperlinks.Add($ws.Cells.Item(3,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/38efd1a5-cded-43d5-9e40-a56e8c5d8d85.328696bf5fa0985729bbcefda211a578ce768408.zh-cn.xlf", "", "", "38efd1a5-cded-43d5-9e40-a56e8c5d8d85.328696bf5fa0985729bbcefda211a578ce768408.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Cells.Item(4,1), "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/5e154c0c-7e94-41fe-895d-735eefa2dd79.md", "", "", "5e154c0c-7e94-41fe-895d-735eefa2dd79.md")
$ws.Hyperlinks.Add($ws.Cells.Item(4,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5e154c0c-7e94-41fe-895d-735eefa2dd79.9c7f514114b497b10602888709ccdf35bc06a240.zh-cn.xlf", "", "", "5e154c0c-7e94-41fe-895d-735eefa2dd79.9c7f514114b497b10602888709ccdf35bc06a240.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Cells.Item(5,1), "https://github.com/OpenLocalizationTest/oltest/blob/7e65c197f1b142ae59144fc72bb407b8868b3994/e2e/e011a692-b796-4695-a3af-0d2a96468c8d.md", "", "", "e011a692-b796-4695-a3af-0d2a96468c8d.md")
$ws.Hyperlinks.Add($ws.Cells.Item(5,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a0d8d1a8dec46161aa203de9fc5d27573f9aab7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e011a692-b796-4695-a3af-0d2a96468c8d.2f245ddf99038b9ed485dc1f9e82a9223c0ed617.zh-cn.xlf", "", "", "e011a692-b796-4695-a3af-0d2a96468c8d.2f245ddf99038b9ed485dc1f9e82a9223c0ed617.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A3:A4").EntireRow.Insert()
$ws.Hyperlinks.Delete()

$ws.Range("A3").Value = "38efd1a5-cded-43d5-9e40-a56e8c5d8d85.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "38efd1a5-cded-43d5-9e40-a56e8c5d8d85.328696bf5fa0985729bbcefda211a578ce768408.de-de.xlf"
$ws.Range("E3").Value = "2016-03-21 10:34:43"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("J3").Value = "Include"

$ws.Range("A4").Value = "5e154c0c-7e94-41fe-895d-735eefa2dd79.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "5e154c0c-7e94-41fe-895d-735eefa2dd79.9c7f514114b497b10602888709ccdf35bc06a240.de-de.xlf"
$ws.Range("E4").Value = "2016-03-21 10:34:43"
$ws.Range("H4").Value = "0001-01-01 00:00:00"
$ws.Range("J4").Value = "Include"

# Row 5 already holds the shifted-down e011a692 values/format, untouched.

$ws.Hyperlinks.Add($ws.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/007b4ece8461d1d54c433d65ab9d1e56c718bffb/e2e/38cdfdee-e750-404e-bd36-1a8a0262a62d.md", "", "", "38cdfdee-e750-404e-bd36-1a8a0262a62d.md")
$ws.Hyperlinks.Add($ws.Cells.Item(2,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/968e17860d0e5ded421d45e193c41646f4a0852d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/38cdfdee-e750-404e-bd36-1a8a0262a62d.e55059b9176280270322352189b79dcaa829d7ae.de-de.xlf", "", "", "38cdfdee-e750-404e-bd36-1a8a0262a62d.e55059b9176280270322352189b79dcaa829d7ae.de-de.xlf")
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/aa94c032958fea58b6e92f78ca8daff42791b31c/e2e/38cdfdee-e750-404e-bd36-1a8a0262a62d.md", "", "", "38cdfdee-e750-404e-bd36-1a8a0262a62d.md")
$ws.Hyperlinks.Add($ws.Cells.Item(2,7), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/519559a0c996ceded01468ceb56d15fae264c3ba/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/38cdfdee-e750-404e-bd36-1a8a0262a62d.e55059b9176280270322352189b79dcaa829d7ae.de-de.xlf", "", "", "38cdfdee-e750-404e-bd36-1a8a0262a62d.e55059b9176280270322352189b79dcaa829d7ae.de-de.xlf")

$ws.Hyperlinks.Add($ws.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/38efd1a5-cded-43d5-9e40-a56e8c5d8d85.md", "", "", "38efd1a5-cded-43d5-9e40-a56e8c5d8d85.md")
$ws.Hyperlinks.Add($ws.Cells.Item(3,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/38efd1a5-cded-43d5-9e40-a56e8c5d8d85.328696bf5fa0985729bbcefda211a578ce768408.de-de.xlf", "", "", "38efd1a5-cded-43d5-9e40-a56e8c5d8d85.328696bf5fa0985729bbcefda211a578ce768408.de-de.xlf")

$ws.Hyperlinks.Add($ws.Cells.Item(4,1), "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/5e154c0c-7e94-41fe-895d-735eefa2dd79.md", "", "", "5e154c0c-7e94-41fe-895d-735eefa2dd79.md")
$ws.Hyperlinks.Add($ws.Cells.Item(4,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5e154c0c-7e94-41fe-895d-735eefa2dd79.9c7f514114b497b10602888709ccdf35bc06a240.de-de.xlf", "", "", "5e154c0c-7e94-41fe-895d-735eefa2dd79.9c7f514114b497b10602888709ccdf35bc06a240.de-de.xlf")

$ws.Hyperlinks.Add($ws.Cells.Item(5,1), "https://github.com/OpenLocalizationTest/oltest/blob/7e65c197f1b142ae59144fc72bb407b8868b3994/e2e/e011a692-b796-4695-a3af-0d2a96468c8d.md", "", "", "e011a692-b796-4695-a3af-0d2a96468c8d.md")
$ws.Hyperlinks.Add($ws.Cells.Item(5,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8979b13b8b276c62fd185304bb039c5b22074548/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e011a692-b796-4695-a3af-0d2a96468c8d.2f245ddf99038b9ed485dc1f9e82a9223c0ed617.de-de.xlf", "", "", "e011a692-b796-4695-a3af-0d2a96468c8d.2f245ddf99038b9ed485dc1f9e82a9223c0ed617.de-de.xlf")
